$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ C = -1.2966612431;         D = -7.33257971608425;  E = 4.73925722988425;  F = 0.952100623564654;  G = 0.998885744847557 }
    3  = @{ C = -0.472246910500001;    D = -6.50816538348425;  E = 5.56367156248425;  F = 0.998885744847557;  G = 0.998885744847557 }
    4  = @{ C = 1.032046008;           D = -5.00387246498425;  E = 7.06796448098425;  F = 0.97842747434348;   G = 0.998885744847557 }
    5  = @{ C = 1.0573510965;          D = -4.97856737648425;  E = 7.09326956948425;  F = 0.976492978845968;  G = 0.998885744847557 }
    6  = @{ C = -2.1141091599;         D = -5.44348155081401;  E = 1.21526323101401;  F = 0.319979598265584;  G = 0.959938794796751 }
    7  = @{ C = -1.3122898082;         D = -4.64166219911401;  E = 2.01708258271401;  F = 0.714136410670466;  G = 0.998885744847557 }
    8  = @{ C = -1.4063230454;         D = -4.73569543631401;  E = 1.92304934551401;  F = 0.6646383228818;    G = 0.998885744847557 }
    9  = @{ C = 0.400117306199999;     D = -2.92925508471401;  E = 3.72948969711401;  F = 0.9941624118913;    G = 0.998885744847557 }
    10 = @{ C = 1.701374947;           D = -5.1561841126464;   E = 8.55893400664639;  F = 0.922919369977848;  G = 0.998885744847557 }
    11 = @{ C = 7.269737673;           D = 0.412178613353602;  E = 14.1272967326464;  F = 0.0345026122542675; G = 0.13801044901707 }
    12 = @{ C = 12.7628485872;         D = 5.9052895275536;    E = 19.6204076468464;  F = 0.0000766259297961813; G = 0.000459755578777088 }
    13 = @{ C = 13.5516761597;         D = 6.6941171000536;    E = 20.4092352193464;  F = 0.0000347500707169601; G = 0.000417000848603521 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
